# cryptos.xlsx refresh - GitHub Actions style data update
#
# Rewrites the Price (column D) and Volume(1h) (column E) figures for the
# crypto list, matches one row re-ranking swap (Chainlink / Wrapped liquid
# staked Ether 2.0 swap places at rows 14-15, Wrapped Ether / Polygon swap
# places at rows 16-17) and replaces the last listed coin (row 51,
# HuobiToken -> RocketPoolETH) with its link/price/volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    # Writes $Text into ($Row, $Col) while forcing it to be stored as
    # text, mirroring the workbook's original inline-string cells.
    # Values such as '271.15' or '2.67' look numeric, and a plain
    # `.Value = ...` assignment would let Excel coerce them into real
    # floating point numbers (and mangle them with binary rounding).
    # Pre-setting the NumberFormat to Text ('@') keeps the literal
    # string; re-applying the 'Normal' style afterwards discards the
    # now-unneeded number-format override so no stray style index is
    # left behind on cells that originally had none.
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '43.769.86'  # D2
Set-TextCell 2 5 '  +0.26%  '  # E2
# Row 3
Set-TextCell 3 4 '2.229.06'  # D3
Set-TextCell 3 5 '  +1.90%  '  # E3
# Row 4
Set-TextCell 4 5 '  +0.01%  '  # E4
# Row 5
Set-TextCell 5 4 '271.15'  # D5
Set-TextCell 5 5 '  +5.47%  '  # E5
# Row 6
Set-TextCell 6 4 '92.91'  # D6
Set-TextCell 6 5 '  +15.87%  '  # E6
# Row 7
Set-TextCell 7 4 '0.625'  # D7
Set-TextCell 7 5 '  +0.44%  '  # E7
# Row 8
Set-TextCell 8 5 '  +0.00%  '  # E8
# Row 9
Set-TextCell 9 4 '0.623'  # D9
Set-TextCell 9 5 '  +5.62%  '  # E9
# Row 10
Set-TextCell 10 4 '46.13'  # D10
Set-TextCell 10 5 '  +8.08%  '  # E10
# Row 11
Set-TextCell 11 4 '0.0973'  # D11
Set-TextCell 11 5 '  +6.13%  '  # E11
# Row 12
Set-TextCell 12 4 '8.34'  # D12
Set-TextCell 12 5 '  +20.34%  '  # E12
# Row 13
Set-TextCell 13 5 '  +1.76%  '  # E13
# Row 14
Set-TextCell 14 2 'Chainlink'  # B14
Set-TextCell 14 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'  # C14
Set-TextCell 14 4 '15.11'  # D14
Set-TextCell 14 5 '  +6.25%  '  # E14
# Row 15
Set-TextCell 15 2 'WrappedliquidstakedEther2.0'  # B15
Set-TextCell 15 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'  # C15
Set-TextCell 15 4 '2.561.29'  # D15
Set-TextCell 15 5 '  +1.60%  '  # E15
# Row 16
Set-TextCell 16 2 'WrappedEther'  # B16
Set-TextCell 16 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'  # C16
Set-TextCell 16 4 '2.236.83'  # D16
Set-TextCell 16 5 '  +3.09%  '  # E16
# Row 17
Set-TextCell 17 2 'Polygon'  # B17
Set-TextCell 17 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'  # C17
Set-TextCell 17 4 '0.801'  # D17
Set-TextCell 17 5 '  +3.79%  '  # E17
# Row 18
Set-TextCell 18 4 '43.738.54'  # D18
Set-TextCell 18 5 '  +0.13%  '  # E18
# Row 19
Set-TextCell 19 4 '0.0000105'  # D19
Set-TextCell 19 5 '  +2.83%  '  # E19
# Row 20
Set-TextCell 20 4 '6.07'  # D20
Set-TextCell 20 5 '  +2.83%  '  # E20
# Row 21
Set-TextCell 21 4 '70.66'  # D21
Set-TextCell 21 5 '  +0.87%  '  # E21
# Row 22
Set-TextCell 22 5 '  -1.20%  '  # E22
# Row 23
Set-TextCell 23 4 '233.82'  # D23
Set-TextCell 23 5 '  +1.77%  '  # E23
# Row 24
Set-TextCell 24 4 '9.14'  # D24
Set-TextCell 24 5 '  +3.34%  '  # E24
# Row 25
Set-TextCell 25 5 '  +0.01%  '  # E25
# Row 26
Set-TextCell 26 4 '11.43'  # D26
Set-TextCell 26 5 '  +8.18%  '  # E26
# Row 27
Set-TextCell 27 5 '  +13.61%  '  # E27
# Row 28
Set-TextCell 28 5 '  +4.92%  '  # E28
# Row 29
Set-TextCell 29 4 '41.08'  # D29
Set-TextCell 29 5 '  +0.03%  '  # E29
# Row 30
Set-TextCell 30 4 '2.24'  # D30
Set-TextCell 30 5 '  +0.12%  '  # E30
# Row 31
Set-TextCell 31 4 '172.29'  # D31
Set-TextCell 31 5 '  -0.33%  '  # E31
# Row 32
Set-TextCell 32 4 '0.0922'  # D32
Set-TextCell 32 5 '  +6.35%  '  # E32
# Row 33
Set-TextCell 33 4 '20.94'  # D33
Set-TextCell 33 5 '  +3.08%  '  # E33
# Row 34
Set-TextCell 34 5 '  +4.94%  '  # E34
# Row 35
Set-TextCell 35 5 '  +2.12%  '  # E35
# Row 36
Set-TextCell 36 5 '  -0.29%  '  # E36
# Row 37
Set-TextCell 37 5 '  -0.04%  '  # E37
# Row 38
Set-TextCell 38 5 '  -2.54%  '  # E38
# Row 39
Set-TextCell 39 4 '3.57'  # D39
Set-TextCell 39 5 '  +25.37%  '  # E39
# Row 40
Set-TextCell 40 4 '12.92'  # D40
Set-TextCell 40 5 '  -0.98%  '  # E40
# Row 41
Set-TextCell 41 4 '0.222'  # D41
Set-TextCell 41 5 '  +12.50%  '  # E41
# Row 42
Set-TextCell 42 4 '2.17'  # D42
Set-TextCell 42 5 '  +3.82%  '  # E42
# Row 43
Set-TextCell 43 4 '63.76'  # D43
Set-TextCell 43 5 '  +2.91%  '  # E43
# Row 44
Set-TextCell 44 5 '  -1.60%  '  # E44
# Row 45
Set-TextCell 45 4 '0.0996'  # D45
Set-TextCell 45 5 '  +1.55%  '  # E45
# Row 46
Set-TextCell 46 4 '8.35'  # D46
Set-TextCell 46 5 '  +1.98%  '  # E46
# Row 47
Set-TextCell 47 4 '100.39'  # D47
Set-TextCell 47 5 '  -0.17%  '  # E47
# Row 48
Set-TextCell 48 5 '  +4.71%  '  # E48
# Row 49
Set-TextCell 49 4 '1.19'  # D49
Set-TextCell 49 5 '  +2.84%  '  # E49
# Row 50
Set-TextCell 50 5 '  +1.33%  '  # E50
# Row 51
Set-TextCell 51 2 'RocketPoolETH'  # B51
Set-TextCell 51 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'  # C51
Set-TextCell 51 4 '2.449.36'  # D51
Set-TextCell 51 5 '  +1.66%  '  # E51

